$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 05:35"

# India (row 14) — updated case counts
$ws.Range("B14").Value = 106750
$ws.Range("C14").Value = 275
$ws.Range("E14").Value = 61138
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 3303

# Australia (row 57) — updated case counts
$ws.Range("B57").Value = 7079
$ws.Range("C57").Value = 11
$ws.Range("D57").Value = 6442
$ws.Range("E57").Value = 537

# Haiti's case count overtakes Malta, Chad (Republica del Chad) and Sierra
# Leona, so it moves up to row 125; the displaced countries shift down one
# row each, keeping their own data intact, down to (but not including)
# Jamaica in row 129.
$ws.Range("A125").Value = "Haiti"
$ws.Range("B125").Value = 596
$ws.Range("C125").Value = 63
$ws.Range("D125").Value = 21
$ws.Range("E125").Value = 553
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 22

$ws.Range("A126").Value = "Malta"
$ws.Range("B126").Value = 569
$ws.Range("C126").Value = 0
$ws.Range("D126").Value = 460
$ws.Range("E126").Value = 103
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 6

$ws.Range("A127").Value = "Republica del Chad"
$ws.Range("B127").Value = 545
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 139
$ws.Range("E127").Value = 350
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 56

$ws.Range("A128").Value = "Sierra Leona"
$ws.Range("B128").Value = 534
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 167
$ws.Range("E128").Value = 334
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 33
